$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Friday row (row 7); Saturday shifts up to row 7
$ws.Rows(7).Delete()

# Clear out stale data cells in rows 2-7 (B:F) before repopulating
$ws.Range("B2:F7").ClearContents()

# Populate the lecture table with the new schedule data
$ws.Range("B2").Value = "GEN0806-sec-Hall 5"
$ws.Range("C2").Value = "GEN0806-Cortney Heller-Hall 5"
$ws.Range("D2").Value = "POW1804-sec-Hall 5"
$ws.Range("E2").Value = "GEN1805-Prof. Virgie Braun II-Hall 5"
$ws.Range("F2").Value = "MEC0811-sec-Hall 5"
$ws.Range("B3").Value = "GEN0807-Jayde Predovic I-Hall 5"
$ws.Range("C3").Value = "CIE1808-Lincoln Predovic-Hall 5"
$ws.Range("D3").Value = "CIE2802-Lincoln Predovic-Hall 5"
$ws.Range("E3").Value = "CIE3804-Lexi Cassin-Hall 5"
$ws.Range("F3").Value = "GEN1809-sec-Hall 5"
$ws.Range("B4").Value = "GEN0810-Carroll Hirthe DVM-Hall 5"
$ws.Range("C4").Value = "CIE4818-sec-Hall 5"
$ws.Range("D4").Value = "GEN0801-Percival Greenholt-Hall 5"
$ws.Range("E4").Value = "CIE1808-Amy Cole-Hall 5"
$ws.Range("F4").Value = "CIE3801-Araceli Hand-Hall 5"
$ws.Range("B5").Value = "MEC0811-Carroll Hirthe DVM-Hall 5"
$ws.Range("C5").Value = "GEN1801-sec-Hall 5"
$ws.Range("D5").Value = "POW1804-Carroll Hirthe DVM-Hall 5"
$ws.Range("E5").Value = "GEN1809-Ofelia O'Conner Jr.-Hall 5"
$ws.Range("F5").Value = "CIE3804-lab-Hall 5"
$ws.Range("B6").Value = "GEN0802-lab-Hall 5"
$ws.Range("C6").Value = "CIE3801-sec-Hall 5"
$ws.Range("D6").Value = "CIE2802-sec-Hall 5"
$ws.Range("E6").Value = "GEN0801-sec-Hall 5"
$ws.Range("F6").Value = "GEN1801-Adah Hyatt-Hall 1`nGEN2810-Miss Edna Schuppe-Hall 2`nCIE1808-sec-Hall 5"
$ws.Range("B7").Value = "CIE1803-lab-Hall 1`nGEN0802-Jayde Predovic I-Hall 2`nPOW1804-Danny Prohaska-Hall 5"
$ws.Range("C7").Value = "CIE4818-Percival Greenholt-Hall 2`nCIE1803-Araceli Hand-Hall 5"
$ws.Range("D7").Value = "GEN0801-Marcus Hegmann-Hall 2"
$ws.Range("E7").Value = "GEN1801-Miss Myriam Huel-Hall 1`nGEN0807-sec-Hall 2`nCIE3801-Lexi Cassin-Hall 5"
$ws.Range("F7").Value = "GEN0809-Marcus Hegmann-Hall 2`nGEN0810-Danny Prohaska-Hall 5"
